$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Cuarentena_HN")

# Update the shared description text in column O (rows 10 to 307).
# All these cells point at the same shared string value, so updating
# the whole column range in one go keeps every occurrence in sync.
$oldText = "Segmentando a la poblacion para poder circular conforme a la terminacion de los digitos de su tarjeta de identidad, pasaporte o carnet de residente para extranjeros, para que puedan abastecerse de insumos básicos,  con horario de 6:00 am a 8:00 pm.  De lunes a domingo, circulando dos digitos por día."
$newText = "Segmentando a la poblacion para poder circular conforme a la terminacion de los digitos de su tarjeta de identidad, pasaporte o carnet de residente para extranjeros, para que puedan abastecerse de insumos básicos,  con horario de 5:00 am a 9:00 pm.  De lunes a domingo."

$rng = $ws.Range("O10:O307")
for ($i = 1; $i -le $rng.Rows.Count; $i++) {
    $cell = $rng.Cells.Item($i, 1)
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}

# Update the visible view: scroll so column F is the top-left visible
# column, and move the active selection to K12.
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("K12").Select()
